$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.208.00"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "3.741.18"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "3.738.14"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "4.365.37"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "3.750.86"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "69.279.68"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "499.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "445.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.25%  "
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").Value = "2.941.34"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("E51").Value = "  -0.60%  "
